$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# ---------------------------------------------------------------------------
# The FilesTab (row 4) Cypher query cell (B4) is rewritten to drop the
# "File Type" and "Breed" output columns (and their coalesce() lines),
# matching the trimmed Bento/ICDC test-script query used for the Files tab.
# ---------------------------------------------------------------------------
$newFilesQuery = "MATCH (f:file)-->(parent)`n" +
  "WITH DISTINCT f, parent`n" +
  "MATCH (f)-[*]->(c:case)<--(demo:demographic)`n" +
  "WHERE demo.breed IN ['Akita']  `n" +
  "OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n" +
  "OPTIONAL MATCH (samp:sample)-->(c)`n" +
  "WITH DISTINCT f, parent, c, demo, diag, s`n" +
  "RETURN  coalesce(f.file_name, '') AS ``File Name``,`n" +
  "        coalesce(labels(parent)[0], '') AS ``Association``,`n" +
  "        coalesce(f.file_description, '') AS ``Description``,`n" +
  "        coalesce(f.file_format, '') AS ``Format``,`n" +
  "        coalesce(f.file_size, '') AS ``Size``,`n" +
  "        coalesce(c.case_id, '') AS ``Case ID``,`n" +
  "        coalesce(diag.disease_term,'') AS Diagnosis , `n" +
  "        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFilesQuery

# The row shrank (two fewer output columns / coalesce lines), so Excel's
# auto row height for the wrapped-text cell drops from 246.5 to 217.5.
$ws.Rows.Item(4).RowHeight = 217.5

# The saved workbook's cursor/viewport moved down to the Files row.
$ws.Range("B4").Select()

$wb.Save()
